$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 3000
$ws.Range("I7").Value = 4000
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 4000
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = -3888
$ws.Range("N7").Value = -2224
$ws.Range("H10").Value = 3000
$ws.Range("J10").Value = 3000
$ws.Range("L10").Value = 3000
$ws.Range("N10").Value = -3586
$ws.Range("H14").Value = 3000
$ws.Range("I14").Value = 4000
$ws.Range("J14").Value = 2000
$ws.Range("K14").Value = 4000
$ws.Range("L14").Value = 2000
$ws.Range("M14").Value = -3809
$ws.Range("N14").Value = -2382
$ws.Range("H33").Value = 561.55554
$ws.Range("I33").Value = 561.55554
$ws.Range("K33").Value = 561.55554
$ws.Range("M33").Value = -332.55554
$ws.Range("H40").Value = 5333.3335
$ws.Range("I40").Value = 5500
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 5500
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -5325
$ws.Range("N40").Value = -5350
$ws.Range("H87").Value = 279000
$ws.Range("J87").Value = 279000
$ws.Range("L87").Value = 279000
$ws.Range("N87").Value = -281496
$ws.Range("H88").Value = 1949.3
$ws.Range("J88").Value = 2055.889
$ws.Range("L88").Value = 2055.889
$ws.Range("N88").Value = -2867.889
$ws.Range("H90").Value = 279000
$ws.Range("J90").Value = 279000
$ws.Range("L90").Value = 837000
$ws.Range("N90").Value = -849480
$ws.Range("H91").Value = 1949.3
$ws.Range("J91").Value = 2055.889
$ws.Range("L91").Value = 2055.889
$ws.Range("N91").Value = -4863.889
$ws.Range("H92").Value = 55556104
$ws.Range("I92").Value = 83333740
$ws.Range("K92").Value = 83333740
$ws.Range("M92").Value = -83332492
$ws.Range("H100").Value = 3744.5557
$ws.Range("I100").Value = 3117.6667
$ws.Range("K100").Value = 3117.6667
$ws.Range("M100").Value = -2576.6667
$ws.Range("H116").Value = 3445
$ws.Range("I116").Value = 3445
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 3445
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -3
$ws.Range("H132").Value = 1772.7333
$ws.Range("I132").Value = 1772.7333
$ws.Range("K132").Value = 5318.199900000001
$ws.Range("M132").Value = -2788.199900000001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 514.1429000000001
$ws.Range("I4").Value = 514.1429000000001
$ws.Range("K4").Value = 514.1429000000001
$ws.Range("M4").Value = -398.1429000000001
$ws.Range("H122").Value = 10833
$ws.Range("I122").Value = 8749.5
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 26248.5
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -23798.5
$ws.Range("N122").Value = -49900

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").ClearContents()
$ws.Range("N117").Value = 0

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 3942.8572
$ws.Range("I32").Value = 2619.8
$ws.Range("K32").Value = 2619.8
$ws.Range("M32").Value = -2303.8
$ws.Range("H103").Value = 41126
$ws.Range("I103").Value = 36647
$ws.Range("K103").Value = 36647
$ws.Range("M103").Value = -35475

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H80").Value = 1789.6
$ws.Range("I80").Value = 1499
$ws.Range("K80").Value = 4497
$ws.Range("M80").Value = -3561
$ws.Range("H83").Value = 1789.6
$ws.Range("I83").Value = 1499
$ws.Range("K83").Value = 13491
$ws.Range("M83").Value = -8811
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = 0
$ws.Range("H97").Value = 1207
$ws.Range("I97").Value = 999.75
$ws.Range("J97").Value = 1483.3334
$ws.Range("K97").Value = 2999.25
$ws.Range("L97").Value = 4450.0002
$ws.Range("M97").Value = -2503.25
$ws.Range("N97").Value = -5442.0002

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H128").Value = 79665.336
$ws.Range("J128").Value = 79665.336
$ws.Range("L128").Value = 79665.336
$ws.Range("N128").Value = -89625.336

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2502.3333
$ws.Range("I7").Value = 2465.6667
$ws.Range("K7").Value = 2465.6667
$ws.Range("M7").Value = -2353.6667
$ws.Range("H22").Value = 453.53125
$ws.Range("I22").Value = 225.9375
$ws.Range("J22").Value = 681.125
$ws.Range("K22").Value = 225.9375
$ws.Range("L22").Value = 681.125
$ws.Range("M22").Value = 69.0625
$ws.Range("N22").Value = -1271.125
$ws.Range("H27").Value = 453.53125
$ws.Range("I27").Value = 225.9375
$ws.Range("J27").Value = 681.125
$ws.Range("K27").Value = 225.9375
$ws.Range("L27").Value = 681.125
$ws.Range("M27").Value = -118.9375
$ws.Range("N27").Value = -895.125
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H100").Value = 2727.2727
$ws.Range("I100").Value = 2710.5
$ws.Range("K100").Value = 2710.5
$ws.Range("M100").Value = -2169.5
$ws.Range("H122").Value = 4863.4287
$ws.Range("I122").Value = 4125
$ws.Range("J122").Value = 5848
$ws.Range("K122").Value = 12375
$ws.Range("L122").Value = 17544
$ws.Range("M122").Value = -9925
$ws.Range("N122").Value = -22444
$ws.Range("H126").Value = 2502.3333
$ws.Range("I126").Value = 2465.6667
$ws.Range("K126").Value = 7397.000100000001
$ws.Range("M126").Value = -4927.000100000001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6730.8
$ws.Range("I126").Value = 1874.5
$ws.Range("K126").Value = 5623.5
$ws.Range("M126").Value = -3153.5
